$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2413793103448276
$ws.Range("C2").Value = 0.4482758620689655
$ws.Range("J2").Value = 0.03448275862068965
$ws.Range("P2").Value = 0.2068965517241379
$ws.Range("S2").Value = 0.06896551724137931
$ws.Range("B3").Value = 0.07692307692307693
$ws.Range("P3").Value = 0.8461538461538461
$ws.Range("S3").Value = 0.07692307692307693
$ws.Range("B6").Value = 0.125
$ws.Range("J6").Value = 0.375
$ws.Range("Q6").Value = 0.375
$ws.Range("S6").Value = 0.125
$ws.Range("B7").Value = 0.2
$ws.Range("J7").Value = 0.4
$ws.Range("Q7").Value = 0.4
$ws.Range("B8").Value = 0.1538461538461539
$ws.Range("D8").Value = 0.03846153846153846
$ws.Range("J8").Value = 0.1538461538461539
$ws.Range("O8").Value = 0.03846153846153846
$ws.Range("Q8").Value = 0.2307692307692308
$ws.Range("S8").Value = 0.3076923076923077
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.2666666666666667
$ws.Range("S9").Value = 0.2
$ws.Range("B10").Value = 0.09933774834437085
$ws.Range("D10").Value = 0.01986754966887417
$ws.Range("F10").Value = 0.03311258278145696
$ws.Range("J10").Value = 0.1390728476821192
$ws.Range("O10").Value = 0.006622516556291391
$ws.Range("Q10").Value = 0.3311258278145696
$ws.Range("R10").Value = 0.09271523178807947
$ws.Range("S10").Value = 0.2781456953642384
$ws.Range("G11").Value = 0.2857142857142857
$ws.Range("J11").Value = 0.1428571428571428
$ws.Range("K11").Value = 0.2857142857142857
$ws.Range("L11").Value = 0.2857142857142857
$ws.Range("G12").Value = 0.5
$ws.Range("J12").Value = 0.5
$ws.Range("J13").Value = 1
$ws.Range("H15").Value = 0.2857142857142857
$ws.Range("J15").Value = 0.5
$ws.Range("S15").Value = 0.2142857142857143
$ws.Range("H16").Value = 0.05
$ws.Range("I16").Value = 0.05
$ws.Range("J16").Value = 0.6
$ws.Range("K16").Value = 0.05
$ws.Range("O16").Value = 0.05
$ws.Range("S16").Value = 0.2
$ws.Range("H17").Value = 0.140625
$ws.Range("I17").Value = 0.109375
$ws.Range("J17").Value = 0.625
$ws.Range("M17").Value = 0.015625
$ws.Range("N17").Value = 0.015625
$ws.Range("O17").Value = 0.046875
$ws.Range("S17").Value = 0.046875
$ws.Range("H18").Value = 0.05
$ws.Range("I18").Value = 0.05
$ws.Range("J18").Value = 0.7
$ws.Range("K18").Value = 0.05
$ws.Range("O18").Value = 0.05
$ws.Range("S18").Value = 0.1
$ws.Range("F19").Value = 0.01388888888888889
$ws.Range("H19").Value = 0.1527777777777778
$ws.Range("I19").Value = 0.08333333333333333
$ws.Range("J19").Value = 0.5694444444444444
$ws.Range("K19").Value = 0.04166666666666666
$ws.Range("O19").Value = 0.06944444444444445
$ws.Range("S19").Value = 0.06944444444444445
